$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:250 down to 173:251
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with its data
$ws.Cells.Item(172, 1).Value = 10
$ws.Cells.Item(172, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(172, 3).Value = "La Araucanía"
$ws.Cells.Item(172, 4).Value = 45016
$ws.Cells.Item(172, 4).NumberFormat = $ws.Cells.Item(173, 4).NumberFormat
$ws.Cells.Item(172, 5).Value = 9
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100104
$ws.Cells.Item(172, 8).Value = "Frutos de pepita"
$ws.Cells.Item(172, 9).Value = 100104003
$ws.Cells.Item(172, 10).Value = "Membrillo"
$ws.Cells.Item(172, 11).Value = "Champion"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 95
$ws.Cells.Item(172, 14).Value = 14000
$ws.Cells.Item(172, 15).Value = 14000
$ws.Cells.Item(172, 16).Value = 14000
$ws.Cells.Item(172, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(172, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(172, 19).Value = 778
$ws.Cells.Item(172, 20).Value = 18
